$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Tabelle24 (M:Q) ---------------------------------------------------
# Row 18 keeps the existing "Google Maps API einbauen" task (now finished
# at 40%) - fill it in before M17 is overwritten so the shared string for
# that text is preserved instead of being reused/renamed in place.
$ws.Range("M18").Value = "Google Maps API einbauen"
$ws.Range("N18").Value = 42690
$ws.Range("O18").Value = 0.59027777777777779
$ws.Range("P18").Value = 0.65972222222222221
$ws.Range("Q18").Value = 0.4

# Row 17 becomes the new "Sprint 2 Planung" task, fully completed.
$ws.Range("M17").Value = "Sprint 2 Planung"
$ws.Range("Q17").Value = 1

# --- Tabelle245 (W:AA) --------------------------------------------------
# New row 18 entry: "Sprint 2 Planung", completed.
$ws.Range("W18").Value = "Sprint 2 Planung"
$ws.Range("X18").Value = 42690
$ws.Range("Y18").Value = 0.59027777777777779
$ws.Range("Z18").Value = 0.65972222222222221
$ws.Range("AA18").Value = 1

# --- Tabelle2 (C:G) ------------------------------------------------------
# New row 20 entry: "Sprint 2 Planung", completed.
$ws.Range("C20").Value = "Sprint 2 Planung"
$ws.Range("D20").Value = 42690
$ws.Range("E20").Value = 0.59027777777777779
$ws.Range("F20").Value = 0.65972222222222221
$ws.Range("G20").Value = 1

# New row 21 entry: "Implementierung Administrative Verwaltungs-Applikation", completed.
$ws.Range("C21").Value = "Implementierung Administrative Verwaltungs-Applikation"
$ws.Range("D21").Value = 42693
$ws.Range("E21").Value = 0.77083333333333337
$ws.Range("F21").Value = 0.875
$ws.Range("G21").Value = 1

# --- View state ----------------------------------------------------------
# Reflects the new selection left after filling in the rows above
# (scroll position itself isn't persisted by this host, but we still set
# it for completeness / forward-compatibility).
$excel.ActiveWindow.ScrollColumn = 14
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("G20:G21").Select() | Out-Null
